# Weekly data update: a new price-report row is inserted right after the
# current row 108 (becoming row 109), pushing all the existing rows
# (109..140) down by one (110..141). We then populate the freshly
# inserted row 109 with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 109; existing rows 109-140 shift down to 110-141.
$ws.Rows(109).Insert()

# Populate the new row 109 with the new week's record. Columns that are
# constant across this data block (A, B, C, E, F, G, H, I, N, Q, R) are
# copied from the template; the varying columns (D, J, K, L, M, O, P) get
# this week's reported figures.
$ws.Range("A109").Value = 11
$ws.Range("B109").Value = "Vega Monumental Concepción"
$ws.Range("C109").Value = "Bíobío"
$ws.Range("D109").Value = 44985
$ws.Range("E109").Value = 8
$ws.Range("F109").Value = 100112001
$ws.Range("G109").Value = "Berenjena"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 100
$ws.Range("K109").Value = 7500
$ws.Range("L109").Value = 8000
$ws.Range("M109").Value = 7750
$ws.Range("N109").Value = "`$/caja 60 unidades"
$ws.Range("O109").Value = "Región de Arica y Parinacota"
$ws.Range("P109").Value = 129
$ws.Range("Q109").Value = 60
$ws.Range("R109").Value = "Hortaliza"
